# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.810.20"
$ws.Range("E2").Value = "  +5.50%  "
$ws.Range("D3").Value = "2.759.68"
$ws.Range("E3").Value = "  +5.49%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "332.55"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("E7").Value = "  +2.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.84"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0861"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.49%  "
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.65"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.54%  "
$ws.Range("D15").Value = "3.190.16"
$ws.Range("E15").Value = "  +5.26%  "
$ws.Range("D16").Value = "2.763.14"
$ws.Range("E16").Value = "  +5.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.887"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.62%  "
$ws.Range("D18").Value = "51.733.06"
$ws.Range("E18").Value = "  +5.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.23"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.46"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.06%  "
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  +3.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.37"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.61"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("E25").Value = "  +5.42%  "
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.09"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.98"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0824"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.03"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.99"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0352"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +10.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "127.11"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.46%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.02"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("E44").Value = "  +7.58%  "
$ws.Range("E45").Value = "  +13.96%  "
$ws.Range("D46").Value = "2.089.02"
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("E47").Value = "  +3.58%  "
$ws.Range("E48").Value = "  +5.39%  "
$ws.Range("E49").Value = "  +6.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.05"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.01"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.04%  "
